# Page object model sheet working = drug , location , QTY
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table Data")

# Update the "To/From (Location)" column (E) values that were placeholder
# numbers / stale locations to the correct facility / pharmacy names.
$ws.Range("E4").Value = "External facility"
$ws.Range("E5").Value = "Final Facility"
$ws.Range("E6").Value = "Pharmacy"
$ws.Range("E7").Value = "Internal Facility 2"
$ws.Range("E13").Value = "WA Nursing Home"

# Leave the cursor where the author left it after making the edits.
$ws.Range("E18").Select()
